# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q3" (before "总计").
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($null, $q3)
$newSheet.Name = "2022-Q1"

# Copy the header / index-column formatting from the "2021-Q3" sheet so the
# new sheet keeps the same visual style (bordered, bold, centered cells).
$q3.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q3.Range("A2:A7").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)   # xlPasteFormats

# Headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row index column (A) is numeric, 0-based.
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("A6").Value = 4
$newSheet.Range("A7").Value = 5

# Columns B-G are stored as text (fund codes keep leading zeros, numeric
# looking figures stay text just like the existing sheets do).
$newSheet.Range("B2:G7").NumberFormat = "@"

$newSheet.Range("B2").Value = "002423"

$newSheet.Range("C2").Value = "华宝兴业标普美国消费(QDII-LOF)美元"
$newSheet.Range("D2").Value = "3.62"
$newSheet.Range("E2").Value = "94.37"
$newSheet.Range("F2").Value = "3.92"
$newSheet.Range("G2").Value = "0.1419"
$newSheet.Range("H2").Value = 5

$newSheet.Range("B3").Value = "162415"
$newSheet.Range("C3").Value = "华宝标普美国消费(QDII-LOF)人民币A"
$newSheet.Range("D3").Value = "3.62"
$newSheet.Range("E3").Value = "94.37"
$newSheet.Range("F3").Value = "3.92"
$newSheet.Range("G3").Value = "0.1419"
$newSheet.Range("H3").Value = 5

$newSheet.Range("B4").Value = "000593"
$newSheet.Range("C4").Value = "易方达标普全球高端消费品指数增强(QDII)-美元现汇"
$newSheet.Range("D4").Value = "1.93"
$newSheet.Range("E4").Value = "92.46"
$newSheet.Range("F4").Value = "6.46"
$newSheet.Range("G4").Value = "0.1247"
$newSheet.Range("H4").Value = 6

$newSheet.Range("B5").Value = "005676"
$newSheet.Range("C5").Value = "易方达标普全球高端消费品指数增强C(QDII) - 人民币"
$newSheet.Range("D5").Value = "1.93"
$newSheet.Range("E5").Value = "92.46"
$newSheet.Range("F5").Value = "6.46"
$newSheet.Range("G5").Value = "0.1247"
$newSheet.Range("H5").Value = 6

$newSheet.Range("B6").Value = "118002"
$newSheet.Range("C6").Value = "易方达标普全球高端消费品指数增强A(QDII) - 人民币"
$newSheet.Range("D6").Value = "1.93"
$newSheet.Range("E6").Value = "92.46"
$newSheet.Range("F6").Value = "6.46"
$newSheet.Range("G6").Value = "0.1247"
$newSheet.Range("H6").Value = 6

$newSheet.Range("B7").Value = "009975"
$newSheet.Range("C7").Value = "华宝标普美国消费(QDII-LOF)人民币C"
$newSheet.Range("D7").Value = "0.61"
$newSheet.Range("E7").Value = "94.37"
$newSheet.Range("F7").Value = "3.92"
$newSheet.Range("G7").Value = "0.0239"
$newSheet.Range("H7").Value = 5

# The "@" text-number-format above left a stray style on B2:G7; strip it
# back out (without touching the already-stored text values/types) by
# pasting in the plain/unstyled format from a known unformatted cell.
$q3.Range("C2").Copy()
$newSheet.Range("B2:G7").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new first data row for the
#    2022-Q1 quarter, pushing the existing rows down and renumbering the
#    index column.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# The freshly inserted row inherits stray formatting from the insert - clear
# it so the new cells match the plain (unstyled) look of the other rows.
$summary.Range("B2:D2").ClearFormats()

# Re-apply the bordered/bold style used by the index column to the new A2.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q1"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.68

# Renumber the index column for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2

Write-Output "done"
